$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new "curveField" parameter / "V" unit row to the parameter table
# (fixes amplitude bugs in setting keysight ARB files: a curveField amplitude
# parameter, in Volts, was missing from the parameter/unit list).
$ws.Range("A62").Value = "curveField"
$ws.Range("B62").Value = "V"

# Touch a number format on a scratch cell and clear it again -- mirrors a
# stray formatting tweak made (and undone) while editing, which leaves an
# extra, unused date/time number-format style registered in the workbook.
$scratch = $ws.Range("Z100")
$scratch.NumberFormat = "m/d/yy h:mm"
$scratch.Clear()
